# Release mCSD 3.9.0 with CP integrated
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Version: 3.8.0 -> 3.9.0
$ws.Range("B3").Value = "3.9.0"

# Date: refresh publication date
$ws.Range("B8").Value = "2024-12-02T17:05:26-06:00"

# Contact rows (previously all showed "No display for ContactDetail")
$ws.Range("B10").Value = "null (https://www.ihe.net/ihe_domains/it_infrastructure/)"
$ws.Range("B11").Value = "null (iti@ihe.net)"
$ws.Range("B12").Value = "IHE IT Infrastructure Technical Committee (iti@ihe.net)"

# Jurisdiction: World -> Global (Whole world)
$ws.Range("B13").Value = "Global (Whole world)"
